$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 267
$endRow = 279

$dates = @(
    "12/03/2020",
    "12/04/2020",
    "12/05/2020",
    "12/06/2020",
    "12/07/2020",
    "12/08/2020",
    "12/09/2020",
    "12/10/2020",
    "12/11/2020",
    "12/12/2020",
    "12/13/2020",
    "12/14/2020",
    "12/15/2020"
)

$sp = @(0.4, 0.39, 0.41, 0.46, 0.4, 0.41, 0.39, 0.4, 0.4, 0.41, 0.45, 0.4, 0.4)
$mogi = @(0.39, 0.38, 0.39, 0.43, 0.39, 0.42, 0.38, 0.39, 0.38, 0.4, 0.43, 0.4, 0.41)

$weekdays = @(
    "Quinta-Feira",
    "Sexta-feira",
    "Sábado",
    "Domingo",
    "Segunda-feira",
    "Terça-feira",
    "Quarta-feira",
    "Quinta-Feira",
    "Sexta-feira",
    "Sábado",
    "Domingo",
    "Segunda-feira",
    "Terça-feira"
)

# Force column A to text format so the date-like strings aren't
# auto-converted to date serial numbers (column D values are not
# number/date-like, so they don't need this).
$colARange = $ws.Range("A$startRow`:A$endRow")
$colARange.NumberFormat = "@"

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $sp[$i]
    $ws.Cells.Item($r, 3).Value = $mogi[$i]
    $ws.Cells.Item($r, 4).Value = $weekdays[$i]
}

# Restore default style (no explicit style index) on column A so the
# new cells match the unstyled data rows already in the sheet.
$colARange.Style = "Normal"
